$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1958456973293768
$ws.Range("C2").Value = 0.5578635014836796
$ws.Range("J2").Value = 0.02077151335311573
$ws.Range("P2").Value = 0.1602373887240356
$ws.Range("S2").Value = 0.06528189910979229
$ws.Range("B3").Value = 0.01522842639593909
$ws.Range("C3").Value = 0.03553299492385787
$ws.Range("J3").Value = 0.03045685279187817
$ws.Range("P3").Value = 0.7208121827411168
$ws.Range("S3").Value = 0.1979695431472081
$ws.Range("P4").Value = 0.7608695652173914
$ws.Range("S4").Value = 0.2391304347826087
$ws.Range("B6").Value = 0.09871244635193133
$ws.Range("D6").Value = 0.0128755364806867
$ws.Range("F6").Value = 0.09012875536480687
$ws.Range("J6").Value = 0.1974248927038627
$ws.Range("O6").Value = 0.01716738197424893
$ws.Range("Q6").Value = 0.1545064377682404
$ws.Range("R6").Value = 0.07725321888412018
$ws.Range("S6").Value = 0.351931330472103
$ws.Range("B7").Value = 0.1098901098901099
$ws.Range("D7").Value = 0.04395604395604396
$ws.Range("F7").Value = 0.08241758241758242
$ws.Range("J7").Value = 0.1208791208791209
$ws.Range("O7").Value = 0.01648351648351648
$ws.Range("Q7").Value = 0.1868131868131868
$ws.Range("R7").Value = 0.07692307692307693
$ws.Range("S7").Value = 0.3626373626373626
$ws.Range("B8").Value = 0.1011494252873563
$ws.Range("D8").Value = 0.01149425287356322
$ws.Range("F8").Value = 0.05977011494252873
$ws.Range("J8").Value = 0.128735632183908
$ws.Range("O8").Value = 0.02298850574712644
$ws.Range("Q8").Value = 0.1839080459770115
$ws.Range("R8").Value = 0.1103448275862069
$ws.Range("S8").Value = 0.3816091954022989
$ws.Range("B9").Value = 0.0893854748603352
$ws.Range("D9").Value = 0.01675977653631285
$ws.Range("E9").Value = 0.0111731843575419
$ws.Range("F9").Value = 0.0893854748603352
$ws.Range("J9").Value = 0.1340782122905028
$ws.Range("O9").Value = 0.0111731843575419
$ws.Range("Q9").Value = 0.1787709497206704
$ws.Range("R9").Value = 0.1229050279329609
$ws.Range("S9").Value = 0.3463687150837989
$ws.Range("B10").Value = 0.1234567901234568
$ws.Range("D10").Value = 0.02083333333333333
$ws.Range("E10").Value = 0.0007716049382716049
$ws.Range("F10").Value = 0.07330246913580248
$ws.Range("J10").Value = 0.1350308641975309
$ws.Range("O10").Value = 0.0162037037037037
$ws.Range("Q10").Value = 0.209104938271605
$ws.Range("R10").Value = 0.09722222222222222
$ws.Range("S10").Value = 0.3240740740740741
$ws.Range("F11").Value = 0.003610108303249098
$ws.Range("G11").Value = 0.1083032490974729
$ws.Range("J11").Value = 0.1010830324909747
$ws.Range("K11").Value = 0.1805054151624549
$ws.Range("L11").Value = 0.592057761732852
$ws.Range("S11").Value = 0.01444043321299639
$ws.Range("G12").Value = 0.7590361445783133
$ws.Range("J12").Value = 0.1927710843373494
$ws.Range("K12").Value = 0.01204819277108434
$ws.Range("L12").Value = 0.02409638554216868
$ws.Range("S12").Value = 0.01204819277108434
$ws.Range("G13").Value = 0.673469387755102
$ws.Range("J13").Value = 0.2857142857142857
$ws.Range("S13").Value = 0.04081632653061224
$ws.Range("F15").Value = 0.02459016393442623
$ws.Range("H15").Value = 0.1885245901639344
$ws.Range("I15").Value = 0.05737704918032787
$ws.Range("J15").Value = 0.3442622950819672
$ws.Range("K15").Value = 0.05327868852459016
$ws.Range("M15").Value = 0.01229508196721311
$ws.Range("O15").Value = 0.1147540983606557
$ws.Range("S15").Value = 0.2049180327868853
$ws.Range("F16").Value = 0.03111111111111111
$ws.Range("H16").Value = 0.1822222222222222
$ws.Range("I16").Value = 0.05777777777777778
$ws.Range("J16").Value = 0.4666666666666667
$ws.Range("K16").Value = 0.09333333333333334
$ws.Range("M16").Value = 0.01777777777777778
$ws.Range("O16").Value = 0.04
$ws.Range("S16").Value = 0.1111111111111111
$ws.Range("F17").Value = 0.01108647450110865
$ws.Range("H17").Value = 0.1995565410199557
$ws.Range("I17").Value = 0.07760532150776053
$ws.Range("J17").Value = 0.4390243902439024
$ws.Range("K17").Value = 0.09090909090909091
$ws.Range("M17").Value = 0.02882483370288248
$ws.Range("N17").Value = 0.002217294900221729
$ws.Range("O17").Value = 0.06873614190687362
$ws.Range("S17").Value = 0.082039911308204
$ws.Range("F18").Value = 0.03524229074889868
$ws.Range("H18").Value = 0.1629955947136564
$ws.Range("I18").Value = 0.1101321585903084
$ws.Range("J18").Value = 0.3964757709251101
$ws.Range("K18").Value = 0.09251101321585903
$ws.Range("M18").Value = 0.013215859030837
$ws.Range("O18").Value = 0.07488986784140969
$ws.Range("S18").Value = 0.1145374449339207
$ws.Range("F19").Value = 0.0091324200913242
$ws.Range("H19").Value = 0.2027397260273973
$ws.Range("I19").Value = 0.08493150684931507
$ws.Range("J19").Value = 0.3963470319634703
$ws.Range("K19").Value = 0.1159817351598174
$ws.Range("M19").Value = 0.02557077625570776
$ws.Range("O19").Value = 0.08127853881278539
$ws.Range("S19").Value = 0.08401826484018265
